$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value = "TestAutomation16oct"
$ws.Range("B8").Value = "TestAutomation16oct"
$ws.Range("C8").Value = "Facility_POC16oct"
$ws.Range("D8").Value = "Facility_POC16oct"
$ws.Range("E8").Value = "Pharmacy_POC16oct"
$ws.Range("F8").Value = "Pharmacy_POC16oct"

$ws.Range("H11").Select()
